{"js": "// Auto-generated Office.js edit script\nconst changes = [\n  {\"find\": \"The rainbow tables for even just a million passwords can be nearly 40GB in size.  The problem we address is how a distributed system could be used to provide fast lookup of matching passwords in this table for some given capture\", \"replace\": \"The rainbow tables for even just a million passwords can be gigabytes in size.  The problem we address is how a distributed system could be used to provide fast lookup of matching passwords in a table for some given capture\"},\n  {\"find\": \"  We chose to focus on WPA versus other wireless encryption techniques (e.g. WEP) because it offers the best encryption that cannot be feasibly defeated via brute force methods.\", \"replace\": \"  We chose to focus on WPA versus other wireless encryption techniques (e.g. WEP) because it offers the best encryption that cannot currently be feasibly defeated via brute force methods.\"},\n  {\"find\": \"Our results show an order of magnitude of 8 times in speed increase for our distributed system versus the serial coWPAtty.  While our data shows that serial coWPAtty could still return a single result in an order of several seconds versus our system returning results in under 1 second the scalability of our system would provide much more usability to offer this as a service.\", \"replace\": \"Our results show an order of magnitude of 8 times in speed increase for our distributed system versus the serial coWPAtty.  Testing data showed that the original serial coWPAtty code could still return a single result in an order of several seconds while our system returning results in under 1 second.  While not a seemingly big difference the scalability of our system would provide much more usability to offer this as a service.\"},\n  {\"find\": \"  The most important effect from our research shows that using generated large rainbow tables of pre\", \"replace\": \"  The most important outcome from our research shows that using large rainbow tables of pre\"},\n  {\"find\": \" easy, fast, and scalable tool for finding weak passwords in encrypted wireless networks.\", \"replace\": \" easy, fast, and scalable method for finding weak passwords in encrypted wireless networks.\"},\n  {\"find\": \"Early advents of these networks provided security of their networks to limit access and protect sensitive data with a \", \"replace\": \"Early advents of these networks provided security to limit access and protect sensitive data with a \"},\n  {\"find\": \"The deficiency of this protocol lead to the development of a new one know\", \"replace\": \"The deficiencies of this protocol lead to the development of a new one know\"},\n  {\"find\": \"There are multiple types of WPA encryption.  Some which use client certificates\", \"replace\": \"There are multiple types of WPA encryption some which use client certificates\"},\n  {\"find\": \"The worker nodes were created by modifying the original coWPAtty code to function in a distributed manner and are described laster.\", \"replace\": \"The worker nodes were created by modifying the original coWPAtty code to function in a distributed manner and are described later.\"},\n  {\"find\": \"Common binaries for all the workers such as the actual worker binary executable is stored here as well.\", \"replace\": \"Common binaries for all the workers such as the actual worker binary executable are stored here as well.\"},\n  {\"find\": \"When a user submits a job via the interface the username as reported by the web server container (Apache Tomcat) will be used for the name of the job owner and as part of the job id.\", \"replace\": \"When a user submits a job via the interface the username, as reported by the web server container (Apache Tomcat), will be used for the name of the job owner and as part of the job id.\"},\n  {\"find\": \"Information such as the start and end time along with the solution found if any is listed in the interface.\", \"replace\": \"Information such as the start and end time along with the solution found, if any, is listed in the interface.\"},\n  {\"find\": \"Details on the various state types is listed in the Worker Node section later.\", \"replace\": \"Details on the various state types are listed in the Worker Node section later.\"},\n  {\"find\": \"The master uses an ssh connection to remotely connect to the worker node specified by the configuration and issues a command to start the worker.  The executable is accessible via a network file system share and an appropriate pathname is used in the command to this binary executable file.\", \"replace\": \"The master uses an SSH connection to remotely connect to the worker node specified by the configuration and issues a command to start the worker.  The executable is accessible via a network file system share, and an appropriate pathname is used in the command to this binary executable file.\"},\n  {\"find\": \"LOADED \u2013 the worker is ready to accept TCP connections, has loaded the rainbow table, and has not yet run any jobs.\", \"replace\": \"LOADED \u2013 the worker is ready to accept jobs, has loaded the rainbow table, and has not yet run any jobs.\"},\n  {\"find\": \"The packet itself consists of multiple field values that are always null terminated in the actual value as well as separated with the special control character\", \"replace\": \"The packet itself consists of multiple field values that are always null terminated as well as separated with the special control character\"},\n  {\"find\": \"Jobid is a uniquely generated id from the master for logging purposes\", \"replace\": \"jobid is a uniquely generated id from the master for logging purposes\"},\n  {\"find\": \"The user request is assigned a job id and added to the queue.\", \"replace\": \"The user\u2019s request is assigned a job id and added to the queue.\"},\n  {\"find\": \"If one of the workers reports back to the master that it found the solution\", \"replace\": \"If one of the workers found the solution\"},\n  {\"find\": \"The master tells all the other workers to stop\", \"replace\": \"The master tells all the workers to stop\"},\n  {\"find\": \"If all of the workers report back that they didn\u2019t find the solution\", \"replace\": \"If none of the workers found the solution\"},\n  {\"find\": \"Either way the master\", \"replace\": \"In either case the master\"},\n  {\"find\": \"Updates the display to show the solution or NO SOLUTION\u201d\", \"replace\": \"Updates the display to show the solution or \u201cNO SOLUTION\u201d\"},\n  {\"find\": \"All nine virtual machines has 1 virtual CPU and 200GB of disk storage.  Access to the cluster from the Internet was limited to ssh and https to the master node only.  All worker nodes had an internal IPv4 network on a private vlan on the host machine only.\", \"replace\": \"All nine virtual machines have 1 virtual CPU and 200GB of disk storage allocated.  Access to the cluster from the Internet was limited to SSH and HTTPS to the master node only.  All worker nodes had an internal IPv4 network on a private VLAN on the host machine only.\"},\n  {\"find\": \"The master communicated with the workers over TCP sockets and through ssh remote commands.\", \"replace\": \"The master communicated with the workers over TCP sockets and through SSH remote commands.\"},\n  {\"find\": \"The master node has Oracle Java \", \"replace\": \"The master node had Oracle Java \"},\n  {\"find\": \"binaries.  The actual rainbow table was hosted on local disk for each node to provide better performance during loading of the node software.\", \"replace\": \"binaries.  The actual rainbow table was hosted on local disk on each node to provide better performance during loading of the node software.\"},\n  {\"find\": \"Acknowledgment of the great work done by Joshua Wright, the original developer of coWPAtty from which our work is based.\", \"replace\": \"We give acknowledgment of the great work done by Joshua Wright, the original developer of coWPAtty, from which our work is based.\"},\n  {\"find\": \"WARNING if you reorder these you may have to fix [#] entires in the paper itself manually!!!\", \"replace\": \"WARNING if you reorder these you may have to fix [#] entries in the paper itself manually!!!\"},\n];\n\nfor (const ch of changes) {\n  const results = context.document.body.search(ch.find, { matchCase: true, matchWildcards: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error('Expected exactly 1 match for: ' + ch.find + ' but found ' + results.items.length);\n  }\n  results.items[0].insertText(ch.replace, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Insert two new paragraphs after the anchor paragraph\nconst anchorResults = context.document.body.search(\"TODO MORE INFO IN SUBSECITONS\", { matchCase: true });\nanchorResults.load('items');\nawait context.sync();\nif (anchorResults.items.length !== 1) {\n  throw new Error('Expected exactly 1 match for anchor paragraph, found ' + anchorResults.items.length);\n}\nlet anchorPara = anchorResults.items[0].paragraphs.getFirst();\nfor (const text of [\"You should create more subsections (header type 3)\", \"Describe the mesting methodology, etc.  Maybe more \u201cTEST ENVIORNMENT\u201d into a subsection under \u201cTesting\u201d\"]) {\n  const newPara = anchorPara.insertParagraph(text, Word.InsertLocation.after);\n  newPara.styleBuiltIn = Word.Style.normal;\n  anchorPara = newPara;\n}\nawait context.sync();", "ps1": "# Auto-generated Word COM (PowerShell-style) edit script\n\nfunction Count-Occurrences($doc, $text) {\n    $count = 0\n    $range = $doc.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $text\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop\n    while ($find.Execute()) {\n        $count++\n        $range.Collapse(0)  # wdCollapseEnd\n    }\n    return $count\n}\n\nfunction Replace-UniqueText($doc, $findText, $replaceText) {\n    $n = Count-Occurrences $doc $findText\n    if ($n -ne 1) {\n        throw \"Expected exactly 1 match for [$findText] but found $n\"\n    }\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 0, $false, $replaceText, 1) | Out-Null  # wdReplaceOne\n}\n\nfunction Find-ParagraphIndexByText($doc, $text) {\n    $idx = 0\n    foreach ($p in $doc.Paragraphs) {\n        $idx++\n        $t = $p.Range.Text.TrimEnd([char]13)\n        if ($t -eq $text) {\n            return $idx\n        }\n    }\n    throw \"Paragraph not found: $text\"\n}\n\n$d = $word.ActiveDocument\n\n# Text replacements\nReplace-UniqueText $d \"The rainbow tables for even just a million passwords can be nearly 40GB in size.  The problem we address is how a distributed system could be used to provide fast lookup of matching passwords in this table for some given capture\" \"The rainbow tables for even just a million passwords can be gigabytes in size.  The problem we address is how a distributed system could be used to provide fast lookup of matching passwords in a table for some given capture\"\nReplace-UniqueText $d \"  We chose to focus on WPA versus other wireless encryption techniques (e.g. WEP) because it offers the best encryption that cannot be feasibly defeated via brute force methods.\" \"  We chose to focus on WPA versus other wireless encryption techniques (e.g. WEP) because it offers the best encryption that cannot currently be feasibly defeated via brute force methods.\"\nReplace-UniqueText $d \"Our results show an order of magnitude of 8 times in speed increase for our distributed system versus the serial coWPAtty.  While our data shows that serial coWPAtty could still return a single result in an order of several seconds versus our system returning results in under 1 second the scalability of our system would provide much more usability to offer this as a service.\" \"Our results show an order of magnitude of 8 times in speed increase for our distributed system versus the serial coWPAtty.  Testing data showed that the original serial coWPAtty code could still return a single result in an order of several seconds while our system returning results in under 1 second.  While not a seemingly big difference the scalability of our system would provide much more usability to offer this as a service.\"\nReplace-UniqueText $d \"  The most important effect from our research shows that using generated large rainbow tables of pre\" \"  The most important outcome from our research shows that using large rainbow tables of pre\"\nReplace-UniqueText $d \" easy, fast, and scalable tool for finding weak passwords in encrypted wireless networks.\" \" easy, fast, and scalable method for finding weak passwords in encrypted wireless networks.\"\nReplace-UniqueText $d \"Early advents of these networks provided security of their networks to limit access and protect sensitive data with a \" \"Early advents of these networks provided security to limit access and protect sensitive data with a \"\nReplace-UniqueText $d \"The deficiency of this protocol lead to the development of a new one know\" \"The deficiencies of this protocol lead to the development of a new one know\"\nReplace-UniqueText $d \"There are multiple types of WPA encryption.  Some which use client certificates\" \"There are multiple types of WPA encryption some which use client certificates\"\nReplace-UniqueText $d \"The worker nodes were created by modifying the original coWPAtty code to function in a distributed manner and are described laster.\" \"The worker nodes were created by modifying the original coWPAtty code to function in a distributed manner and are described later.\"\nReplace-UniqueText $d \"Common binaries for all the workers such as the actual worker binary executable is stored here as well.\" \"Common binaries for all the workers such as the actual worker binary executable are stored here as well.\"\nReplace-UniqueText $d \"When a user submits a job via the interface the username as reported by the web server container (Apache Tomcat) will be used for the name of the job owner and as part of the job id.\" \"When a user submits a job via the interface the username, as reported by the web server container (Apache Tomcat), will be used for the name of the job owner and as part of the job id.\"\nReplace-UniqueText $d \"Information such as the start and end time along with the solution found if any is listed in the interface.\" \"Information such as the start and end time along with the solution found, if any, is listed in the interface.\"\nReplace-UniqueText $d \"Details on the various state types is listed in the Worker Node section later.\" \"Details on the various state types are listed in the Worker Node section later.\"\nReplace-UniqueText $d \"The master uses an ssh connection to remotely connect to the worker node specified by the configuration and issues a command to start the worker.  The executable is accessible via a network file system share and an appropriate pathname is used in the command to this binary executable file.\" \"The master uses an SSH connection to remotely connect to the worker node specified by the configuration and issues a command to start the worker.  The executable is accessible via a network file system share, and an appropriate pathname is used in the command to this binary executable file.\"\nReplace-UniqueText $d \"LOADED \u2013 the worker is ready to accept TCP connections, has loaded the rainbow table, and has not yet run any jobs.\" \"LOADED \u2013 the worker is ready to accept jobs, has loaded the rainbow table, and has not yet run any jobs.\"\nReplace-UniqueText $d \"The packet itself consists of multiple field values that are always null terminated in the actual value as well as separated with the special control character\" \"The packet itself consists of multiple field values that are always null terminated as well as separated with the special control character\"\nReplace-UniqueText $d \"Jobid is a uniquely generated id from the master for logging purposes\" \"jobid is a uniquely generated id from the master for logging purposes\"\nReplace-UniqueText $d \"The user request is assigned a job id and added to the queue.\" \"The user\u2019s request is assigned a job id and added to the queue.\"\nReplace-UniqueText $d \"If one of the workers reports back to the master that it found the solution\" \"If one of the workers found the solution\"\nReplace-UniqueText $d \"The master tells all the other workers to stop\" \"The master tells all the workers to stop\"\nReplace-UniqueText $d \"If all of the workers report back that they didn\u2019t find the solution\" \"If none of the workers found the solution\"\nReplace-UniqueText $d \"Either way the master\" \"In either case the master\"\nReplace-UniqueText $d \"Updates the display to show the solution or NO SOLUTION\u201d\" \"Updates the display to show the solution or \u201cNO SOLUTION\u201d\"\nReplace-UniqueText $d \"All nine virtual machines has 1 virtual CPU and 200GB of disk storage.  Access to the cluster from the Internet was limited to ssh and https to the master node only.  All worker nodes had an internal IPv4 network on a private vlan on the host machine only.\" \"All nine virtual machines have 1 virtual CPU and 200GB of disk storage allocated.  Access to the cluster from the Internet was limited to SSH and HTTPS to the master node only.  All worker nodes had an internal IPv4 network on a private VLAN on the host machine only.\"\nReplace-UniqueText $d \"The master communicated with the workers over TCP sockets and through ssh remote commands.\" \"The master communicated with the workers over TCP sockets and through SSH remote commands.\"\nReplace-UniqueText $d \"The master node has Oracle Java \" \"The master node had Oracle Java \"\nReplace-UniqueText $d \"binaries.  The actual rainbow table was hosted on local disk for each node to provide better performance during loading of the node software.\" \"binaries.  The actual rainbow table was hosted on local disk on each node to provide better performance during loading of the node software.\"\nReplace-UniqueText $d \"Acknowledgment of the great work done by Joshua Wright, the original developer of coWPAtty from which our work is based.\" \"We give acknowledgment of the great work done by Joshua Wright, the original developer of coWPAtty, from which our work is based.\"\nReplace-UniqueText $d \"WARNING if you reorder these you may have to fix [#] entires in the paper itself manually!!!\" \"WARNING if you reorder these you may have to fix [#] entries in the paper itself manually!!!\"\n\n# Insert two new paragraphs after the anchor paragraph\n$anchorIdx = Find-ParagraphIndexByText $d \"TODO MORE INFO IN SUBSECITONS\"\n$anchorPara = $d.Paragraphs.Item($anchorIdx)\n$anchorPara.Range.InsertParagraphAfter()\n$anchorIdx++\n$newPara = $d.Paragraphs.Item($anchorIdx)\n$newPara.Range.Text = \"You should create more subsections (header type 3)\"\n$newPara.Style = \"Normal\"\n$anchorPara = $newPara\n$anchorPara.Range.InsertParagraphAfter()\n$anchorIdx++\n$newPara = $d.Paragraphs.Item($anchorIdx)\n$newPara.Range.Text = \"Describe the mesting methodology, etc.  Maybe more \u201cTEST ENVIORNMENT\u201d into a subsection under \u201cTesting\u201d\"\n$newPara.Style = \"Normal\"\n$anchorPara = $newPara\n\nWrite-Output \"Edit complete\""}
